$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Find the "Work on coding files:" bullet list (numId=1) paragraph
#    that currently reads "Integrating a database into code to store
#    user input" and change its text.
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Integrating a database into code to store user input`r") {
        $target = $p
        break
    }
}

$target.Range.Text = "Organized the files for the project and mapped out how it was going to be programmed"

# ------------------------------------------------------------------
# 2. Right after it, insert two new bullets (same numId=1/ilvl=0
#    list formatting, inherited automatically from the paragraph we
#    split):
#       "Researched and assisted with implementing user interface "
#       "Integrating a database into code to store user input"
#    (the latter is the original sentence, now relocated here).
# ------------------------------------------------------------------
$idx = $target.Index
$d.Paragraphs($idx).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs($idx + 1).Range.Text = "Researched and assisted with implementing user interface "

$d.Paragraphs($idx + 1).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs($idx + 2).Range.Text = "Integrating a database into code to store user input"

# ------------------------------------------------------------------
# 3. In Arianna Rodriguez's section, insert a new bullet
#    "Brainstormed project ideas" right before the existing
#    "Work on project documentation:" bullet (numId=2/ilvl=0).
# ------------------------------------------------------------------
$docBullet = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Work on project documentation:`r") {
        $docBullet = $p
        break
    }
}

$docBullet.Range.InsertParagraphBefore() | Out-Null
$docBullet.Range.Text = "Brainstormed project ideas"

# ------------------------------------------------------------------
# 4. Append a brand-new, completely empty paragraph at the very end
#    of the document body (after the last "Film and edit final
#    YouTube video" bullet, before the sectPr).
# ------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>")
